$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the data range as text so numeric-looking strings
# (e.g. "27.725.28", "1.000") are preserved verbatim instead of
# being auto-converted to numbers by the COM Value setter.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '27.725.28'
$ws.Range('E2').Value = '  -0.62%  '
$ws.Range('D3').Value = '1.895.23'
$ws.Range('E3').Value = '  +1.26%  '
$ws.Range('E4').Value = '  -0.87%  '
$ws.Range('D5').Value = '313.12'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('D7').Value = '0.4844'
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('D8').Value = '0.3796'
$ws.Range('E8').Value = '  -0.69%  '
$ws.Range('D9').Value = '0.07338'
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('D10').Value = '0.9155'
$ws.Range('E10').Value = '  -2.57%  '
$ws.Range('D11').Value = '20.54'
$ws.Range('E11').Value = '  -2.62%  '
$ws.Range('D12').Value = '0.07690'
$ws.Range('D13').Value = '1.866.18'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('D14').Value = '5.469'
$ws.Range('E14').Value = '  -0.41%  '
$ws.Range('D15').Value = '6.598'
$ws.Range('E15').Value = '  +0.02%  '
$ws.Range('D16').Value = '90.94'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('E17').Value = '  -0.86%  '
$ws.Range('D18').Value = '0.000008808'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('E19').Value = '  -1.00%  '
$ws.Range('D20').Value = '27.748.97'
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').Value = '14.46'
$ws.Range('E21').Value = '  -2.37%  '
$ws.Range('D22').Value = '5.117'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').Value = '2.105.02'
$ws.Range('E23').Value = '  -1.14%  '
$ws.Range('D24').Value = '10.76'
$ws.Range('E24').Value = '  -0.50%  '
$ws.Range('D25').Value = '1.903'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('D26').Value = '153.92'
$ws.Range('E26').Value = '  -1.49%  '
$ws.Range('D27').Value = '18.37'
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('D28').Value = '2.133'
$ws.Range('E28').Value = '  +3.98%  '
$ws.Range('D29').Value = '115.82'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  -1.73%  '
$ws.Range('D31').Value = '0.08915'
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('D32').Value = '3.151'
$ws.Range('E32').Value = '  -5.37%  '
$ws.Range('D33').Value = '1.227'
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('D34').Value = '0.7646'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('E36').Value = '  -0.34%  '
$ws.Range('D37').Value = '2.527'
$ws.Range('E37').Value = '  -7.46%  '
$ws.Range('D38').Value = '1.094'
$ws.Range('E38').Value = '  -3.59%  '
$ws.Range('D39').Value = '0.05266'
$ws.Range('E39').Value = '  -2.31%  '
$ws.Range('D40').Value = '0.5470'
$ws.Range('E40').Value = '  -3.22%  '
$ws.Range('E41').Value = '  -0.11%  '
$ws.Range('D42').Value = '6.920'
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '8.469'
$ws.Range('E43').Value = '  -1.08%  '
$ws.Range('B44').Value = 'Algorand'
$ws.Range('C44').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D44').Value = '0.1518'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').Value = '110.05'
$ws.Range('E45').Value = '  +4.95%  '
$ws.Range('D46').Value = '10.57'
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('E47').Value = '  -2.42%  '
$ws.Range('D48').Value = '1.000'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D50').Value = '67.40'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').Value = '0.06053'
$ws.Range('E51').Value = '  -0.92%  '

# Restore the default (unstyled) cell style so the text-format
# override above does not leave a stray style index behind.
$ws.Range("B2:E51").Style = "Normal"
